$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value for "Right" answers (row 11, column B)
$ws.Range("B11").Value = 5

# Update total marks for "Right" answers (row 12, column B)
$ws.Range("B12").Value = 100

# Update the correct/total marks summary text (row 12, column E)
$ws.Range("E12").Value = "100/140"
